# Fruta / hortaliza, semanal
# Insert two new weekly price records into the Granada (Vega Modelo de Temuco)
# data table, shifting the existing rows down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at row 122 (pushes old rows 122-134 down to 123-135) ---
$ws.Rows("122:122").Insert()

$ws.Cells.Item(122, 1).Value = 10
$ws.Cells.Item(122, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(122, 3).Value = "La Araucanía"
$ws.Cells.Item(122, 4).Value = 44748
$ws.Cells.Item(122, 5).Value = 9
$ws.Cells.Item(122, 6).Value = "Fruta"
$ws.Cells.Item(122, 7).Value = 100104
$ws.Cells.Item(122, 8).Value = "Frutos de pepita"
$ws.Cells.Item(122, 9).Value = 100104001
$ws.Cells.Item(122, 10).Value = "Granada"
$ws.Cells.Item(122, 11).Value = "Wonderfull"
$ws.Cells.Item(122, 12).Value = "Primera"
$ws.Cells.Item(122, 13).Value = 200
$ws.Cells.Item(122, 14).Value = 13000
$ws.Cells.Item(122, 15).Value = 13000
$ws.Cells.Item(122, 16).Value = 13000
$ws.Cells.Item(122, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(122, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(122, 19).Value = 1300
$ws.Cells.Item(122, 20).Value = 10

# --- Insert second new row at row 133 (pushes old rows 132-134, now at 133-135, down to 134-136) ---
$ws.Rows("133:133").Insert()

$ws.Cells.Item(133, 1).Value = 10
$ws.Cells.Item(133, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(133, 3).Value = "La Araucanía"
$ws.Cells.Item(133, 4).Value = 44747
$ws.Cells.Item(133, 5).Value = 9
$ws.Cells.Item(133, 6).Value = "Fruta"
$ws.Cells.Item(133, 7).Value = 100104
$ws.Cells.Item(133, 8).Value = "Frutos de pepita"
$ws.Cells.Item(133, 9).Value = 100104001
$ws.Cells.Item(133, 10).Value = "Granada"
$ws.Cells.Item(133, 11).Value = "Wonderfull"
$ws.Cells.Item(133, 12).Value = "Primera"
$ws.Cells.Item(133, 13).Value = 50
$ws.Cells.Item(133, 14).Value = 13000
$ws.Cells.Item(133, 15).Value = 13000
$ws.Cells.Item(133, 16).Value = 13000
$ws.Cells.Item(133, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(133, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(133, 19).Value = 1300
$ws.Cells.Item(133, 20).Value = 10
